$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for columns I and J, matching the style of the existing header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# I0 / IF values for rows 2-56 (data rows)
$iValues = @(5,7,5,5,7,5,5,5,6,6,5,8,5,5,6,6,8,5,2,8,5,8,6,6,9,7,4,8,9,7,7,7,4,6,8,8,4,9,2,9,6,9,6,7,6,9,6,8,9,6,6,6,6,7,8)
$jValues = @(5,7,5,5,7,6,6,6,6,6,5,8,6,6,7,6,8,5,2,8,6,8,6,6,9,8,4,8,9,8,9,8,4,6,8,8,4,9,3,9,6,9,7,7,7,9,6,8,9,6,7,6,6,7,8)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}

Write-Output "I0/IF columns added"
